# Apply updated probability values to the "Starting_State" matrix sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$changes = @{
    2 = @{ "B" = 0.2096774193548387; "C" = 0.5282258064516129; "J" = 0.008064516129032258; "P" = 0.1451612903225807; "S" = 0.1088709677419355 }
    3 = @{ "B" = 0.007462686567164179; "C" = 0.02985074626865672; "J" = 0.03731343283582089; "P" = 0.7686567164179104; "S" = 0.1567164179104478 }
    4 = @{ "J" = 0.025; "P" = 0.675; "S" = 0.3 }
    6 = @{ "B" = 0.04; "D" = 0.01142857142857143; "F" = 0.04571428571428571; "J" = 0.2514285714285714; "O" = 0.01142857142857143; "Q" = 0.2514285714285714; "R" = 0.06857142857142857; "S" = 0.32 }
    7 = @{ "B" = 0.106280193236715; "D" = 0.03381642512077294; "F" = 0.04347826086956522; "J" = 0.1304347826086956; "O" = 0.02898550724637681; "Q" = 0.178743961352657; "R" = 0.07246376811594203; "S" = 0.4057971014492754 }
    8 = @{ "B" = 0.07650273224043716; "D" = 0.01366120218579235; "E" = 0.00273224043715847; "F" = 0.06830601092896176; "J" = 0.09016393442622951; "O" = 0.01639344262295082; "Q" = 0.1639344262295082; "R" = 0.09836065573770492; "S" = 0.4699453551912569 }
    9 = @{ "B" = 0.09947643979057591; "D" = 0.01047120418848168; "F" = 0.07329842931937172; "J" = 0.1099476439790576; "O" = 0.02094240837696335; "Q" = 0.1884816753926702; "R" = 0.07853403141361257; "S" = 0.418848167539267 }
    10 = @{ "B" = 0.1028880866425993; "D" = 0.02256317689530686; "E" = 0.0009025270758122744; "F" = 0.05776173285198556; "J" = 0.1137184115523466; "O" = 0.01263537906137184; "Q" = 0.2310469314079422; "R" = 0.08212996389891697; "S" = 0.3763537906137184 }
    11 = @{ "G" = 0.1347517730496454; "J" = 0.07801418439716312; "K" = 0.1702127659574468; "L" = 0.6063829787234043; "S" = 0.01063829787234043 }
    12 = @{ "G" = 0.7624309392265194; "J" = 0.1767955801104972; "L" = 0.03314917127071823; "S" = 0.02762430939226519 }
    13 = @{ "G" = 0.7708333333333334; "J" = 0.1875; "S" = 0.04166666666666666 }
    15 = @{ "F" = 0.004854368932038835; "H" = 0.1213592233009709; "I" = 0.09223300970873786; "J" = 0.383495145631068; "K" = 0.05825242718446602; "M" = 0.01456310679611651; "O" = 0.0825242718446602; "S" = 0.2427184466019418 }
    16 = @{ "F" = 0.00625; "H" = 0.18125; "I" = 0.075; "J" = 0.475; "K" = 0.08125; "M" = 0.03125; "O" = 0.06875000000000001; "S" = 0.08125 }
    17 = @{ "F" = 0.02325581395348837; "H" = 0.1488372093023256; "I" = 0.1046511627906977; "J" = 0.4116279069767442; "K" = 0.09302325581395349; "M" = 0.02093023255813953; "O" = 0.08604651162790698; "S" = 0.1116279069767442 }
    18 = @{ "F" = 0.01183431952662722; "H" = 0.1834319526627219; "I" = 0.08284023668639054; "J" = 0.378698224852071; "K" = 0.1420118343195266; "M" = 0.01775147928994083; "O" = 0.08284023668639054; "S" = 0.1005917159763314 }
    19 = @{ "F" = 0.01347708894878706; "H" = 0.1949685534591195; "I" = 0.09344115004492363; "J" = 0.3647798742138365; "K" = 0.1293800539083558; "M" = 0.02515723270440252; "N" = 0.0008984725965858042; "O" = 0.0664869721473495; "S" = 0.1114106019766397 }
}

foreach ($row in $changes.Keys) {
    $rowChanges = $changes[$row]
    foreach ($col in $rowChanges.Keys) {
        $ws.Range("$col$row").Value = $rowChanges[$col]
    }
}
